# Apply the regenerated CAM/BOM/CPL data to the active workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated placement data for row 29 (was Mid X=143.764, Mid Y=-102.616, Rotation=-90)
$ws.Range("B29").Value = 143.256
$ws.Range("C29").Value = -103.124
$ws.Range("E29").Value = 0

# Updated placement data for row 30 (was Mid X=142.748, Mid Y=-102.616, Rotation=90)
$ws.Range("B30").Value = 143.256
$ws.Range("C30").Value = -101.854
$ws.Range("E30").Value = 0

# Move the active selection from E34 to E31, matching the regenerated sheet view.
$ws.Range("E31").Select()
